# Add four new glucose (Hexose) fragment rules to the rules table.
# These rows mirror the layout/format of the existing "[M+H-Hexose-H2O]+"
# row (row 151), each referencing D151 plus an existing mass-diff row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (151) down into the
# four new rows (152:155), so the new rows keep the same cell styles
# (fill/font) as the rest of the data block.
$ws.Range("A151:G151").Copy()
$ws.Range("A152:G155").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 152: [M+H-Hexose-H2O-CH4]+
$ws.Range("A152").Value = "[M+H-Hexose-H2O-CH4]+"
$ws.Range("B152").Value = 1
$ws.Range("C152").Value = 1
$ws.Range("D152").Formula = "=D151+D91-1.0073"
$ws.Range("E152").Value = 174
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0.5

# Row 153: [M+H-Hexose-H2O-CH3OH]+
$ws.Range("A153").Value = "[M+H-Hexose-H2O-CH3OH]+"
$ws.Range("B153").Value = 1
$ws.Range("C153").Value = 1
$ws.Range("D153").Formula = "=D151+D94-1.0073"
$ws.Range("E153").Value = 175
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0.5

# Row 154: [M+H-Hexose-H2O-C2H4O]+
$ws.Range("A154").Value = "[M+H-Hexose-H2O-C2H4O]+"
$ws.Range("B154").Value = 1
$ws.Range("C154").Value = 1
$ws.Range("D154").Formula = "=D151+D81-1.0073"
$ws.Range("E154").Value = 176
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 0.5

# Row 155: [M+H-Hexose-H2O-C2H4O2]+
$ws.Range("A155").Value = "[M+H-Hexose-H2O-C2H4O2]+"
$ws.Range("B155").Value = 1
$ws.Range("C155").Value = 1
$ws.Range("D155").Formula = "=D151+D110-1.0073"
$ws.Range("E155").Value = 177
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0.5

# Match the final selection/active cell recorded in the saved workbook.
$ws.Range("F155").Select()
